# Updates market-price derived columns (H-N) on several leve-profit rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR worksheets,
# reflecting refreshed market data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Cells.Item(100, 8).Value = 9839056
$ws.Cells.Item(100, 9).Value = 11150464
$ws.Cells.Item(100, 10).Value = 3500
$ws.Cells.Item(100, 11).Value = 11150464
$ws.Cells.Item(100, 12).Value = 3500
$ws.Cells.Item(100, 13).Value = -11149923
$ws.Cells.Item(100, 14).Value = -4582

# Row 112
$ws.Cells.Item(112, 8).Value = 10910275
$ws.Cells.Item(112, 10).Value = 12988256
$ws.Cells.Item(112, 12).Value = 38964768
$ws.Cells.Item(112, 14).Value = -38966984

# Row 133
$ws.Cells.Item(133, 8).Value = 11896.667
$ws.Cells.Item(133, 10).Value = 11896.667
$ws.Cells.Item(133, 12).Value = 11896.667
$ws.Cells.Item(133, 14).Value = -22016.667

# Row 137
$ws.Cells.Item(137, 8).Value = 23810464
$ws.Cells.Item(137, 9).Value = 27778416
$ws.Cells.Item(137, 10).Value = 2753
$ws.Cells.Item(137, 11).Value = 83335248
$ws.Cells.Item(137, 12).Value = 8259
$ws.Cells.Item(137, 13).Value = -83332698
$ws.Cells.Item(137, 14).Value = -13359

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 2417.64
$ws.Cells.Item(61, 9).Value = 1551.5
$ws.Cells.Item(61, 10).Value = 5882.2
$ws.Cells.Item(61, 11).Value = 1551.5
$ws.Cells.Item(61, 12).Value = 5882.2
$ws.Cells.Item(61, 13).Value = -1339.5
$ws.Cells.Item(61, 14).Value = -6306.2

# Row 122
$ws.Cells.Item(122, 8).Value = 1972.7059
$ws.Cells.Item(122, 9).Value = 1868.5
$ws.Cells.Item(122, 10).Value = 2222.8
$ws.Cells.Item(122, 11).Value = 5605.5
$ws.Cells.Item(122, 12).Value = 6668.400000000001
$ws.Cells.Item(122, 13).Value = -3155.5
$ws.Cells.Item(122, 14).Value = -11568.4

# Row 132
$ws.Cells.Item(132, 8).Value = 2838.0454
$ws.Cells.Item(132, 9).Value = 2478.0571
$ws.Cells.Item(132, 10).Value = 4238
$ws.Cells.Item(132, 11).Value = 7434.1713
$ws.Cells.Item(132, 12).Value = 12714
$ws.Cells.Item(132, 13).Value = -4904.1713
$ws.Cells.Item(132, 14).Value = -17774

# Row 136
$ws.Cells.Item(136, 8).Value = 2417.64
$ws.Cells.Item(136, 9).Value = 1551.5
$ws.Cells.Item(136, 10).Value = 5882.2
$ws.Cells.Item(136, 11).Value = 4654.5
$ws.Cells.Item(136, 12).Value = 17646.6
$ws.Cells.Item(136, 13).Value = -2104.5
$ws.Cells.Item(136, 14).Value = -22746.6

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 15626840
$ws.Cells.Item(134, 9).Value = 18183230
$ws.Cells.Item(134, 10).Value = 4458.4443
$ws.Cells.Item(134, 11).Value = 54549690
$ws.Cells.Item(134, 12).Value = 13375.3329
$ws.Cells.Item(134, 13).Value = -54547155
$ws.Cells.Item(134, 14).Value = -18445.3329

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 85601.836
$ws.Cells.Item(16, 9).Value = 251005.5
$ws.Cells.Item(16, 10).Value = 2900
$ws.Cells.Item(16, 11).Value = 251005.5
$ws.Cells.Item(16, 12).Value = 2900
$ws.Cells.Item(16, 13).Value = -250718.5
$ws.Cells.Item(16, 14).Value = -3474

# Row 58
$ws.Cells.Item(58, 8).Value = 1822.5385
$ws.Cells.Item(58, 9).Value = 1160.5
$ws.Cells.Item(58, 10).Value = 3146.6155
$ws.Cells.Item(58, 11).Value = 1160.5
$ws.Cells.Item(58, 12).Value = 3146.6155
$ws.Cells.Item(58, 13).Value = -957.5
$ws.Cells.Item(58, 14).Value = -3552.6155

# Row 113
$ws.Cells.Item(113, 8).Value = 85601.836
$ws.Cells.Item(113, 9).Value = 251005.5
$ws.Cells.Item(113, 10).Value = 2900
$ws.Cells.Item(113, 11).Value = 251005.5
$ws.Cells.Item(113, 12).Value = 2900
$ws.Cells.Item(113, 13).Value = -248835.5
$ws.Cells.Item(113, 14).Value = -7240

# Row 132
$ws.Cells.Item(132, 8).Value = 2265.7551
$ws.Cells.Item(132, 9).Value = 1811.4147
$ws.Cells.Item(132, 11).Value = 5434.2441
$ws.Cells.Item(132, 13).Value = -2904.2441

# Row 136
$ws.Cells.Item(136, 8).Value = 1822.5385
$ws.Cells.Item(136, 9).Value = 1160.5
$ws.Cells.Item(136, 10).Value = 3146.6155
$ws.Cells.Item(136, 11).Value = 3481.5
$ws.Cells.Item(136, 12).Value = 9439.8465
$ws.Cells.Item(136, 13).Value = -931.5
$ws.Cells.Item(136, 14).Value = -14539.8465

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Cells.Item(12, 8).Value = 25.392857
$ws.Cells.Item(12, 9).Value = 15.0625
$ws.Cells.Item(12, 10).Value = 39.166668
$ws.Cells.Item(12, 11).Value = 45.1875
$ws.Cells.Item(12, 12).Value = 117.500004
$ws.Cells.Item(12, 13).Value = 127.8125
$ws.Cells.Item(12, 14).Value = -463.500004

# Row 59
$ws.Cells.Item(59, 8).Value = 2599.6
$ws.Cells.Item(59, 9).Value = 1000
$ws.Cells.Item(59, 10).Value = 2999.5
$ws.Cells.Item(59, 11).Value = 3000
$ws.Cells.Item(59, 12).Value = 8998.5
$ws.Cells.Item(59, 13).Value = -2460
$ws.Cells.Item(59, 14).Value = -10078.5

# Row 114
$ws.Cells.Item(114, 8).Value = 1536.5294
$ws.Cells.Item(114, 9).Value = 1066.6
$ws.Cells.Item(114, 11).Value = 3199.8
$ws.Cells.Item(114, 13).Value = 54.20000000000027

# Row 131
$ws.Cells.Item(131, 8).Value = 5954028.5
$ws.Cells.Item(131, 10).Value = 6946287
$ws.Cells.Item(131, 12).Value = 20838861
$ws.Cells.Item(131, 14).Value = -20848941

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 2224223.5
$ws.Cells.Item(122, 9).Value = 3705036.8
$ws.Cells.Item(122, 10).Value = 3004
$ws.Cells.Item(122, 11).Value = 11115110.4
$ws.Cells.Item(122, 12).Value = 9012
$ws.Cells.Item(122, 13).Value = -11112660.4
$ws.Cells.Item(122, 14).Value = -13912

# Row 123
$ws.Cells.Item(123, 8).Value = 10183.789
$ws.Cells.Item(123, 10).Value = 10183.789
$ws.Cells.Item(123, 12).Value = 10183.789
$ws.Cells.Item(123, 14).Value = -15083.789

$ws = $wb.Worksheets.Item("LTW")
# Row 24
$ws.Cells.Item(24, 8).Value = 3000
$ws.Cells.Item(24, 10).Value = 3000
$ws.Cells.Item(24, 12).Value = 3000
$ws.Cells.Item(24, 14).Value = -3686

# Row 122
$ws.Cells.Item(122, 8).Value = 3085.5356
$ws.Cells.Item(122, 9).Value = 2453.4614
$ws.Cells.Item(122, 11).Value = 7360.3842
$ws.Cells.Item(122, 13).Value = -4910.3842

# Row 136
$ws.Cells.Item(136, 8).Value = 4261.8647
$ws.Cells.Item(136, 9).Value = 2534.0312
$ws.Cells.Item(136, 10).Value = 15320
$ws.Cells.Item(136, 11).Value = 7602.0936
$ws.Cells.Item(136, 12).Value = 45960
$ws.Cells.Item(136, 13).Value = -5052.0936
$ws.Cells.Item(136, 14).Value = -51060

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Cells.Item(100, 8).Value = 3400
$ws.Cells.Item(100, 9).Value = 3400
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 6800
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 13).Value = -6259
$ws.Cells.Item(100, 14).ClearContents()

# Row 113
$ws.Cells.Item(113, 8).Value = 583.25
$ws.Cells.Item(113, 9).Value = 437
$ws.Cells.Item(113, 10).Value = 729.5
$ws.Cells.Item(113, 11).Value = 1311
$ws.Cells.Item(113, 12).Value = 2188.5
$ws.Cells.Item(113, 13).Value = 859
$ws.Cells.Item(113, 14).Value = -6528.5

# Row 123
$ws.Cells.Item(123, 8).Value = 34884.31
$ws.Cells.Item(123, 10).Value = 34884.31
$ws.Cells.Item(123, 12).Value = 34884.31
$ws.Cells.Item(123, 14).Value = -44684.31

# Row 132
$ws.Cells.Item(132, 8).Value = 2043.1075
$ws.Cells.Item(132, 9).Value = 2067.0278
$ws.Cells.Item(132, 10).Value = 1961.0952
$ws.Cells.Item(132, 11).Value = 6201.0834
$ws.Cells.Item(132, 12).Value = 5883.2856
$ws.Cells.Item(132, 13).Value = -3671.0834
$ws.Cells.Item(132, 14).Value = -10943.2856
